$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row removals -----------------------------------------------------
# Two data rows were dropped from the table entirely: the original
# "RM 232" row (row 26) and the original "SC 92" row (row 28). Deleting
# bottom-to-top keeps the row numbers for the earlier delete valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# From here on row numbers refer to the POST-delete (final) layout,
# i.e. the same numbering the target sheet ends up with (A1:F33).

# --- Newly-imputed numeric values (previously blank "NA" cells) -------
$ws.Range("E2").Value = -7.2
$ws.Range("E12").Value = -5.3
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("C30").Value = 11.4
$ws.Range("E31").Value = -8.1
$ws.Range("E33").Value = -10.7

# --- Values newly marked missing (previously numeric, now blank) ------
# A lone apostrophe is how Excel enters an empty, left-quoted text value;
# resetting the style afterwards clears the quote-prefix flag it sets so
# the cell ends up styled just like any other default cell, matching the
# blank "NA" markers already present elsewhere in the sheet.
function Set-Blank($addr) {
  $ws.Range($addr).Value = "'"
  $ws.Range($addr).Style = "Normal"
}

Set-Blank "E6"
Set-Blank "E14"
Set-Blank "E22"
Set-Blank "E23"
Set-Blank "C32"
